# GDE-9497: Initial Commit for Quick Party of PT Health
# Adds a new "PT_Health_SYND" row (rowid 6) to the UAT Deal Scenarios masterlist.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (rowid): numeric 6, formatted like the rows above it (quote-prefixed,
#     text-like style) but additionally left-aligned -----------------------------------
$ws.Range("A7").Value = 6
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)          # xlPasteFormats
$ws.Range("A7").HorizontalAlignment = -4131  # xlLeft

# --- Column B (UAT_Deal_Scenario_Name) -------------------------------------------------
$ws.Range("B7").Value = "PT_Health_SYND"

# --- Column C (Path): reuse the same value/hyperlink-style formatting as C2:C6 --------
$ws.Range("C7").Value = "\DataSet\NewUATDeals_DataSet\"
$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)          # xlPasteFormats

# --- Column D (Filename) ----------------------------------------------------------------
$ws.Range("D7").Value = "Deal_PT_Health_Syndicated.xlsx"

# Clear the marching-ants clipboard state left behind by the copy operations above.
$excel.CutCopyMode = 0

# The user's cursor ends up on B7 after entering the new scenario name.
$ws.Range("B7").Select()
